# Reorders the "Periodo Mora" (E16:E22) labels so the most recent period
# (2408) is listed first and the oldest (2402) last, and keeps each
# "Valor Mora" (F column) amount attached to its original period value
# (2408 -> 36400, 2402 -> 52000), matching the new database ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2408"
$ws.Range("E17").Value = "2407"
$ws.Range("E18").Value = "2406"
$ws.Range("E19").Value = "2405"
$ws.Range("E20").Value = "2404"
$ws.Range("E21").Value = "2403"
$ws.Range("E22").Value = "2402"

$ws.Range("F16").Value = 36400
$ws.Range("F22").Value = 52000
